$wb = $excel.ActiveWorkbook

# Update status text "Ready for handoff" -> "In Translation" on all sheets
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# Adjust column widths to match new (narrower) "Status" content.
# Target stored width is 13.4101845877511 character-units; Excel's
# ColumnWidth setter snaps to a whole-pixel grid on save, so feed it the
# input that lands on the closest reachable grid point (13.333333...).
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
